# Update the "Scopo" (purpose) paragraph text in the document.
$d = $word.ActiveDocument

$find = "sarà quello di creare un veicolo controllabile tramite bluetooth utilizzando un’applicazione sul telefono."
$replace = "Tankino è di realizzare la base di un carro armato con dei motori elettrici. Questi motori dovranno essere controllati da un arduino in modo da poter muovere il carro armato a piacere. Inoltre l’arduino riceverà dei dati da un dispositivo mobile tramite bluetooth in modo da poter controllare il carro armato a distanza. Per poter controllare il carro armato dal dispositivo mobile è inoltre richiesto di sviluppare un applicazione da cui sia possibile controllare il carro armato con dei comandi semplici e intuitivi. "

$range = $d.Content
$range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
